# "intial mobile base test"
#
# Adds two new login-test-data rows to the "loginTestData" sheet, then
# updates the workbook/sheet view state (selected ranges + active sheet)
# to match where the author was working when the file was saved.

$wb = $excel.ActiveWorkbook

# --- Append two new test rows to loginTestData (rows 8-9) -----------------
$ws1 = $wb.Worksheets.Item("loginTestData")

$ws1.Cells.Item(8, 1).Value = "Nho4DGqEoxtXBi1"
$ws1.Cells.Item(8, 2).Value = "RdPS8xmcPF*."
$ws1.Cells.Item(8, 3).Value = "success"

$ws1.Cells.Item(9, 1).Value = "eagermanipulation"
$ws1.Cells.Item(9, 2).Value = "abi3u1nkXd*."
$ws1.Cells.Item(9, 3).Value = "success"

# --- Restore each sheet's last-used selection ------------------------------
$ws3 = $wb.Worksheets.Item("googleForgetUserNamePassword")
$ws3.Range("A28").Select()

$ws1.Range("C7").Select()

# --- loginTestData becomes the active (visible) sheet on open -------------
$ws1.Activate()
